# Refresh cryptos snapshot: updated prices/volumes and three name/link/price
# row-pairs that were re-sorted (rows 14-15, 31-32, 46-48).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.647.53'
$ws.Range('E2').Value = '  +1.93%  '
$ws.Range('D3').Value = '1.902.77'
$ws.Range('E4').Value = '  +0.52%  '
$ws.Range('D5').Value = "'244.92"
$ws.Range('E5').Value = '  +5.24%  '
$ws.Range('D6').Value = "'0.635"
$ws.Range('E6').Value = '  +2.64%  '
$ws.Range('E7').Value = '  +0.46%  '
$ws.Range('D8').Value = "'42.66"
$ws.Range('E8').Value = '  +1.72%  '
$ws.Range('E9').Value = '  +2.88%  '
$ws.Range('E10').Value = '  +1.99%  '
$ws.Range('E11').Value = '  +1.06%  '
$ws.Range('D12').Value = '2.178.18'
$ws.Range('E12').Value = '  +3.06%  '
$ws.Range('D13').Value = "'12.54"
$ws.Range('E13').Value = '  +9.65%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.933.45'
$ws.Range('E14').Value = '  +4.63%  '
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').Value = "'0.692"
$ws.Range('E15').Value = '  +2.80%  '
$ws.Range('D16').Value = "'4.81"
$ws.Range('E16').Value = '  +2.85%  '
$ws.Range('D17').Value = '35.608.10'
$ws.Range('E17').Value = '  +1.80%  '
$ws.Range('D18').Value = "'72.37"
$ws.Range('E18').Value = '  +3.36%  '
$ws.Range('D19').Value = '0.0₃0812'
$ws.Range('E19').Value = '  +2.51%  '
$ws.Range('D20').Value = "'245.32"
$ws.Range('E20').Value = '  +1.88%  '
$ws.Range('D21').Value = "'12.48"
$ws.Range('E21').Value = '  +2.04%  '
$ws.Range('D22').Value = "'4.92"
$ws.Range('E22').Value = '  +3.33%  '
$ws.Range('E23').Value = '  +0.43%  '
$ws.Range('E24').Value = '  +2.41%  '
$ws.Range('D25').Value = "'171.33"
$ws.Range('E25').Value = '  -0.72%  '
$ws.Range('D26').Value = "'2.15"
$ws.Range('E26').Value = '  +31.49%  '
$ws.Range('D27').Value = "'8.42"
$ws.Range('E27').Value = '  +7.69%  '
$ws.Range('D28').Value = "'18.04"
$ws.Range('E28').Value = '  +3.02%  '
$ws.Range('E29').Value = '  +1.78%  '
$ws.Range('D30').Value = "'0.963"
$ws.Range('E30').Value = '  +28.63%  '
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').Value = "'0.0568"
$ws.Range('E31').Value = '  +2.64%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = "'4.10"
$ws.Range('E32').Value = '  +3.48%  '
$ws.Range('E33').Value = '  +0.56%  '
$ws.Range('E34').Value = '  +4.42%  '
$ws.Range('E35').Value = '  +7.69%  '
$ws.Range('E36').Value = '  +4.55%  '
$ws.Range('E37').Value = '  +9.20%  '
$ws.Range('E38').Value = '  +3.64%  '
$ws.Range('E39').Value = '  +4.95%  '
$ws.Range('D40').Value = "'91.55"
$ws.Range('E40').Value = '  +1.89%  '
$ws.Range('D41').Value = '1.364.14'
$ws.Range('E41').Value = '  +1.08%  '
$ws.Range('D42').Value = "'15.43"
$ws.Range('E42').Value = '  +5.61%  '
$ws.Range('E43').Value = '  +12.34%  '
$ws.Range('D44').Value = "'13.16"
$ws.Range('E44').Value = '  +43.29%  '
$ws.Range('E45').Value = '  +3.56%  '
$ws.Range('B46').Value = 'MultiversX'
$ws.Range('C46').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D46').Value = "'46.47"
$ws.Range('E46').Value = '  +36.94%  '
$ws.Range('B47').Value = 'HuobiToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D47').Value = "'2.42"
$ws.Range('E47').Value = '  +0.18%  '
$ws.Range('B48').Value = 'FraxShare'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D48').Value = "'6.73"
$ws.Range('E48').Value = '  +6.18%  '
$ws.Range('E49').Value = '  +0.48%  '
$ws.Range('D50').Value = '2.088.50'
$ws.Range('E50').Value = '  +2.76%  '
$ws.Range('D51').Value = "'3.53"
$ws.Range('E51').Value = '  +3.79%  '
